$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column O (percent_moisture) rows 2-37: replace the computed values with a
# placeholder 0.00001 and clear the 2-decimal-place number formatting back
# to the default "Normal" style (so the cells no longer carry a custom
# numFmt).
$range = $ws.Range("O2:O37")
$range.Style = "Normal"
$range.Value = 0.00001

# Update the active selection to Q24 as recorded by the saved workbook view.
$ws.Range("Q24").Select()
